# Update cryptocurrency price (column D) and volume/% change (column E) data
# for rows 2-51 on the active worksheet, per the upstream GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values scraped for this run: row => (Price, Volume(1h))
$updates = @(
    @{Row=2; D='62.418.88'; E='  +0.20%  '}
    @{Row=3; D='3.033.39'; E='  -0.09%  '}
    @{Row=4; D='0.999'; E='  -0.05%  '}
    @{Row=5; D='542.55'; E='  +1.34%  '}
    @{Row=6; D='134.24'; E='  +0.90%  '}
    @{Row=7; D='1.00'; E='  +0.10%  '}
    @{Row=8; D='3.028.27'; E='  -0.17%  '}
    @{Row=9; D='0.494'; E='  +0.94%  '}
    @{Row=10; D='6.12'; E='  -0.91%  '}
    @{Row=11; D='0.148'; E='  -3.04%  '}
    @{Row=12; D='0.448'; E='  +0.05%  '}
    @{Row=13; D='0.0000223'; E='  +1.59%  '}
    @{Row=14; D='34.42'; E='  +1.83%  '}
    @{Row=15; D='3.517.73'; E='  -0.30%  '}
    @{Row=16; D='62.403.61'; E='  +0.18%  '}
    @{Row=17; D='3.029.89'; E='  -0.03%  '}
    @{Row=18; D='0.108'; E='  -3.33%  '}
    @{Row=19; D='6.64'; E='  +1.19%  '}
    @{Row=20; D='480.39'; E='  +3.67%  '}
    @{Row=21; D='13.30'; E='  +0.40%  '}
    @{Row=22; D='0.675'; E='  -1.51%  '}
    @{Row=23; D='7.07'; E='  +2.06%  '}
    @{Row=24; D='80.87'; E='  +3.81%  '}
    @{Row=25; D='12.17'; E='  +1.77%  '}
    @{Row=26; D='1.00'; E='  +0.19%  '}
    @{Row=27; D='2.71'; E='  +1.60%  '}
    @{Row=28; D='7.81'; E='  +1.39%  '}
    @{Row=29; D='0.998'; E='  -0.23%  '}
    @{Row=30; D='1.94'; E='  +4.56%  '}
    @{Row=31; D='25.75'; E='  +0.03%  '}
    @{Row=32; D='1.13'; E='  -0.87%  '}
    @{Row=33; D='5.70'; E='  +5.52%  '}
    @{Row=34; D='2.37'; E='  +5.55%  '}
    @{Row=35; D='55.11'; E='  -5.07%  '}
    @{Row=36; D='5.89'; E='  +0.12%  '}
    @{Row=37; D='464.51'; E='  +1.42%  '}
    @{Row=38; D='3.162.98'; E='  -0.58%  '}
    @{Row=39; D='0.0804'; E='  +1.92%  '}
    @{Row=40; D='0.0389'; E='  +0.33%  '}
    @{Row=41; D='0.119'; E='  +2.22%  '}
    @{Row=42; D='8.12'; E='  +0.99%  '}
    @{Row=43; D='2.49'; E='  +0.78%  '}
    @{Row=44; D='26.62'; E='  +7.42%  '}
    @{Row=45; D='1.00'; E='  -0.10%  '}
    @{Row=46; D='0.246'; E='  -0.02%  '}
    @{Row=47; D='0.109'; E='  +0.92%  '}
    @{Row=48; D='1.98'; E='  +1.45%  '}
    @{Row=49; D='0.0₃0504'; E='  -1.21%  '}
    @{Row=50; D='114.15'; E='  -5.78%  '}
    @{Row=51; D='1.28'; E='  +3.42%  '}
)

# The Price/Volume columns are stored as plain text (values such as '1.00' or
# '0.999' must stay literal strings, not be coerced into numbers), so force the
# Text number format on the two columns before writing, then restore each cell's
# original style afterwards so no visible formatting changes are introduced.
$dRange = $ws.Range("D2:D51")
$eRange = $ws.Range("E2:E51")
$dOrigStyle = $dRange.Style
$eOrigStyle = $eRange.Style
$dRange.NumberFormat = "@"
$eRange.NumberFormat = "@"

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$dRange.Style = $dOrigStyle
$eRange.Style = $eOrigStyle

Write-Host "Updated $($updates.Count) rows of crypto price/volume data."
